$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037478998811781
$ws.Range("D2").Value = 1.044485507862041
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.051857989883588
$ws.Range("I2").Value = 1.041410640329165
$ws.Range("J2").Value = 1.042581647992885
$ws.Range("K2").Value = 1.047256424613482
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.05460832837209
$ws.Range("N2").Value = 1.044062233518821
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038304265240848
$ws.Range("D3").Value = 1.045135586391539
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.052686663511202
$ws.Range("I3").Value = 1.041620469822508
$ws.Range("J3").Value = 1.043051934727278
$ws.Range("K3").Value = 1.047718171311714
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.055249691056189
$ws.Range("N3").Value = 1.044533188114319
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03883861150471
$ws.Range("D4").Value = 1.045556503135498
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.053223614746326
$ws.Range("I4").Value = 1.041755044488878
$ws.Range("J4").Value = 1.043355912865405
$ws.Range("K4").Value = 1.048016538242527
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.055664789058558
$ws.Range("N4").Value = 1.04483759793627
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.0390633313892
$ws.Range("D5").Value = 1.045733520042528
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.053449525287485
$ws.Range("I5").Value = 1.041811331832452
$ws.Range("J5").Value = 1.043483625342824
$ws.Range("K5").Value = 1.048141871300885
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.055839317061681
$ws.Range("N5").Value = 1.044965491780056
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039101067534252
$ws.Range("D6").Value = 1.045763245634163
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.053487466941913
$ws.Range("I6").Value = 1.041820765822659
$ws.Range("J6").Value = 1.043505064107299
$ws.Range("K6").Value = 1.048162909348377
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.055868622255546
$ws.Range("N6").Value = 1.044986960990036
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038841613907611
$ws.Range("D7").Value = 1.045558868196698
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.053226632683874
$ws.Range("I7").Value = 1.041755797735147
$ws.Range("J7").Value = 1.043357619681376
$ws.Range("K7").Value = 1.048018213345086
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.055667121030511
$ws.Range("N7").Value = 1.044839307176116
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03775782935741
$ws.Range("D8").Value = 1.044705147495145
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.052137888844922
$ws.Range("I8").Value = 1.041481800983225
$ws.Range("J8").Value = 1.042740650955044
$ws.Range("K8").Value = 1.047412559040942
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.054825059598462
$ws.Range("N8").Value = 1.044221462283425
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035850764087873
$ws.Range("D9").Value = 1.043202945253058
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.050225163720742
$ws.Range("I9").Value = 1.040989835413988
$ws.Range("J9").Value = 1.041651008233089
$ws.Range("K9").Value = 1.046342203753529
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.05334202493069
$ws.Range("N9").Value = 1.04313027214382
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034581296458568
$ws.Range("D10").Value = 1.042203033417966
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.048954011310369
$ws.Range("I10").Value = 1.040655754543713
$ws.Range("J10").Value = 1.040922994177115
$ws.Range("K10").Value = 1.045626614704874
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.052353956345516
$ws.Range("N10").Value = 1.042401224224382
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034032075206929
$ws.Range("D11").Value = 1.041770451217198
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.048404559847381
$ws.Range("I11").Value = 1.040509655919595
$ws.Range("J11").Value = 1.040607395020054
$ws.Range("K11").Value = 1.045316293748312
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.051926278679542
$ws.Range("N11").Value = 1.042085176880319
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03382814162499
$ws.Range("D12").Value = 1.041609830685118
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.048200616109566
$ws.Range("I12").Value = 1.040455172899255
$ws.Range("J12").Value = 1.040490113910807
$ws.Range("K12").Value = 1.045200957952566
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.05176744614342
$ws.Range("N12").Value = 1.041967729218444
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033871882803725
$ws.Range("D13").Value = 1.041644281602534
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.048244356051391
$ws.Range("I13").Value = 1.040466869430943
$ws.Range("J13").Value = 1.040515273510015
$ws.Range("K13").Value = 1.045225700964252
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.051801515060236
$ws.Range("N13").Value = 1.041992924547169
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034015216522969
$ws.Range("D14").Value = 1.04175717304486
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.048387698784767
$ws.Range("I14").Value = 1.040505156727551
$ws.Range("J14").Value = 1.040597701612062
$ws.Range("K14").Value = 1.04530676145116
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.05191314899327
$ws.Range("N14").Value = 1.042075469706576
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034103538697817
$ws.Range("D15").Value = 1.04182673715864
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.04847603650886
$ws.Range("I15").Value = 1.040528718263579
$ws.Range("J15").Value = 1.040648481228704
$ws.Range("K15").Value = 1.045356696416799
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.051981933834834
$ws.Range("N15").Value = 1.042126321436097
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034617756426755
$ws.Range("D16").Value = 1.042231750759177
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.048990497105839
$ws.Range("I16").Value = 1.040665420367185
$ws.Range("J16").Value = 1.040943931881257
$ws.Range("K16").Value = 1.045647199970722
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.052382343462494
$ws.Range("N16").Value = 1.042422191662466
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034940437635525
$ws.Range("D17").Value = 1.042485909717261
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.049313464579203
$ws.Range("I17").Value = 1.040750785157935
$ws.Range("J17").Value = 1.04112916359046
$ws.Range("K17").Value = 1.045829301257982
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.052633554471626
$ws.Range("N17").Value = 1.042607686421943
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035128697042792
$ws.Range("D18").Value = 1.042634193486341
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.049501939113357
$ws.Range("I18").Value = 1.040800438063123
$ws.Range("J18").Value = 1.04123717088592
$ws.Range("K18").Value = 1.045935472794102
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.052780097181792
$ws.Range("N18").Value = 1.042715847100151
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035192896207499
$ws.Range("D19").Value = 1.042684760667932
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.049566219788093
$ws.Range("I19").Value = 1.040817344831852
$ws.Range("J19").Value = 1.041273992571042
$ws.Range("K19").Value = 1.045971666862171
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.05283006705359
$ws.Range("N19").Value = 1.042752721076289
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034905812338472
$ws.Range("D20").Value = 1.042458637011402
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.049278803585475
$ws.Range("I20").Value = 1.040741640688718
$ws.Range("J20").Value = 1.041109293616089
$ws.Range("K20").Value = 1.045809768165292
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.052606600285066
$ws.Range("N20").Value = 1.042587788229928
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033973006325006
$ws.Range("D21").Value = 1.041723927665331
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.048345483849495
$ws.Range("I21").Value = 1.040493888015294
$ws.Range("J21").Value = 1.040573430058567
$ws.Range("K21").Value = 1.045282893046553
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.051880274852305
$ws.Range("N21").Value = 1.042051163684691
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033386929441674
$ws.Range("D22").Value = 1.041262333335418
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.047759520114002
$ws.Range("I22").Value = 1.04033686971944
$ws.Range("J22").Value = 1.040236202253747
$ws.Range("K22").Value = 1.044951228458217
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.051423756379001
$ws.Range("N22").Value = 1.041713456977707
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033697580007359
$ws.Range("D23").Value = 1.041506999792424
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.048070069243902
$ws.Range("I23").Value = 1.040420225972674
$ws.Range("J23").Value = 1.040415001995086
$ws.Range("K23").Value = 1.045127087369106
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.051665750607867
$ws.Range("N23").Value = 1.041892510635187
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034921457881165
$ws.Range("D24").Value = 1.042470960258669
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.04929446510923
$ws.Range("I24").Value = 1.040745773110102
$ws.Range("J24").Value = 1.041118272111197
$ws.Range("K24").Value = 1.045818594468473
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.052618779673794
$ws.Range("N24").Value = 1.04259677947553
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036343456459404
$ws.Range("D25").Value = 1.043591033699923
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.050718952523519
$ws.Range("I25").Value = 1.041118099586906
$ws.Range("J25").Value = 1.041932991712921
$ws.Range("K25").Value = 1.046619277611801
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.053725322275387
$ws.Range("N25").Value = 1.043412656072538
